$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 2919
$ws.Cells.Item(10, 6).Value = 6812
$ws.Cells.Item(12, 6).Value = 38
$ws.Cells.Item(13, 6).Value = 341
$ws.Cells.Item(15, 6).Value = 1478
$ws.Cells.Item(16, 6).Value = 1101
$ws.Cells.Item(17, 6).Value = 2212
$ws.Cells.Item(18, 6).Value = 1459
$ws.Cells.Item(20, 6).Value = 101
$ws.Cells.Item(21, 6).Value = 1096
$ws.Cells.Item(22, 6).Value = 107
$ws.Cells.Item(23, 6).Value = 168
$ws.Cells.Item(24, 6).Value = 326
$ws.Cells.Item(25, 6).Value = 1672
$ws.Cells.Item(26, 6).Value = 1660
$ws.Cells.Item(28, 6).Value = 1026
$ws.Cells.Item(29, 6).Value = 32
$ws.Cells.Item(30, 6).Value = 1652
$ws.Cells.Item(31, 6).Value = 1192
$ws.Cells.Item(32, 6).Value = 132
$ws.Cells.Item(36, 6).Value = 402
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(38, 6).Value = 2436
$ws.Cells.Item(39, 6).Value = 2694
$ws.Cells.Item(41, 6).Value = 183
$ws.Cells.Item(45, 6).Value = 313
$ws.Cells.Item(47, 6).Value = 163
$ws.Cells.Item(48, 6).Value = 140
$ws.Cells.Item(49, 6).Value = 410

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 181
$ws.Cells.Item(7, 6).Value = 156
$ws.Cells.Item(14, 6).Value = 57
$ws.Cells.Item(20, 6).Value = 31
$ws.Cells.Item(22, 6).Value = 337
$ws.Cells.Item(23, 6).Value = 460
$ws.Cells.Item(32, 6).Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 1724
$ws.Cells.Item(7, 6).Value = 1658
$ws.Cells.Item(8, 6).Value = 1843
$ws.Cells.Item(9, 6).Value = 2700
$ws.Cells.Item(10, 6).Value = 989
$ws.Cells.Item(11, 6).Value = 891
$ws.Cells.Item(13, 6).Value = 240
$ws.Cells.Item(14, 6).Value = 1374
$ws.Cells.Item(15, 6).Value = 7245

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 1724
$ws.Cells.Item(6, 6).Value = 2919
$ws.Cells.Item(8, 6).Value = 1658
$ws.Cells.Item(10, 6).Value = 2700
$ws.Cells.Item(11, 6).Value = 6812
$ws.Cells.Item(12, 6).Value = 989
$ws.Cells.Item(13, 6).Value = 891
$ws.Cells.Item(15, 6).Value = 341
$ws.Cells.Item(16, 6).Value = 156
$ws.Cells.Item(17, 6).Value = 240
$ws.Cells.Item(18, 6).Value = 1374
$ws.Cells.Item(20, 6).Value = 2212
$ws.Cells.Item(21, 6).Value = 1459
$ws.Cells.Item(23, 6).Value = 101
$ws.Cells.Item(24, 6).Value = 1096
$ws.Cells.Item(25, 6).Value = 107
$ws.Cells.Item(26, 6).Value = 326
$ws.Cells.Item(27, 6).Value = 57
$ws.Cells.Item(28, 6).Value = 1672
$ws.Cells.Item(29, 6).Value = 1026
$ws.Cells.Item(31, 6).Value = 32
$ws.Cells.Item(32, 6).Value = 1652
$ws.Cells.Item(33, 6).Value = 1192
$ws.Cells.Item(35, 6).Value = 337
$ws.Cells.Item(36, 6).Value = 460
$ws.Cells.Item(37, 6).Value = 402
$ws.Cells.Item(39, 6).Value = 2436
$ws.Cells.Item(40, 6).Value = 2694
$ws.Cells.Item(42, 6).Value = 183
$ws.Cells.Item(44, 6).Value = 313
$ws.Cells.Item(46, 6).Value = 163
$ws.Cells.Item(48, 6).Value = 410
